$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.049.33"
$ws.Range("E2").Value = "  -2.84%  "
$ws.Range("D3").Value = "2.368.71"
$ws.Range("E3").Value = "  -1.82%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "502.29"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.09%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "130.88"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.00%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.57%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.548"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.51%  "
$ws.Range("D9").Value = "2.374.51"
$ws.Range("E9").Value = "  -3.12%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0972"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.03%  "
$ws.Range("E11").Value = "  +0.45%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.327"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.63%  "
$ws.Range("E13").Value = "  +0.04%  "
$ws.Range("D14").Value = "2.790.85"
$ws.Range("E14").Value = "  -1.88%  "
$ws.Range("D15").Value = "55.981.92"
$ws.Range("E15").Value = "  -2.76%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.44"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.45%  "
$ws.Range("E17").Value = "  -2.08%  "
$ws.Range("D18").Value = "2.293.52"
$ws.Range("E18").Value = "  -6.02%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.02"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.13%  "
$ws.Range("E20").Value = "  -3.25%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "307.04"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.06%  "
$ws.Range("E22").Value = "  -2.02%  "
$ws.Range("E23").Value = "  +0.28%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.14"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.95%  "
$ws.Range("E25").Value = "  +0.98%  "
$ws.Range("E26").Value = "  -3.91%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.149"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.31%  "
$ws.Range("E28").Value = "  -4.28%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "172.82"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.33%  "
$ws.Range("E30").Value = "  -3.06%  "
$ws.Range("E31").Value = "  -3.26%  "
$ws.Range("E32").Value = "  +0.23%  "
$ws.Range("E33").Value = "  -5.13%  "
$ws.Range("E34").Value = "  -7.74%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.87%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.61"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.14%  "
$ws.Range("E37").Value = "  -4.37%  "
$ws.Range("E38").Value = "  -1.97%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.01"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.82%  "
$ws.Range("E40").Value = "  -3.05%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.41"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.54%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "131.08"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.06%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.36"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.05%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.79"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -5.85%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.565"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.37%  "
$ws.Range("E46").Value = "  -0.69%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "243.74"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -6.83%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0480"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.77%  "
$ws.Range("E49").Value = "  -3.02%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "17.01"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.02%  "
$ws.Range("E51").Value = "  -2.84%  "
